$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("hpi")
$ws.Range("D2").Value = 'Regurgitation relief suggests that the pain may be related to reflux rather than a structural issue like esophageal stricture.'
$ws.Range("C3").Value = 'Hoarse voice reported is absent.'
$ws.Range("D3").Value = 'A hoarse voice can indicate laryngeal involvement or irritation, which is less common in isolated esophageal stricture cases.'
$ws.Range("C4").Value = 'Cough reported is absent.'
$ws.Range("D4").Value = 'The absence of cough may suggest that there is no significant irritation or obstruction affecting the airway, which can be associated with esophageal issues.'
$ws.Range("A5").Value = 'Antacids no longer providing relief is present.'
$ws.Range("B5").Value = 'The ineffectiveness of antacids may indicate a structural issue like esophageal stricture, as it suggests that the underlying problem is not acid-related.'
$ws.Range("C5").Value = 'Nausea and/or vomiting is absent.'
$ws.Range("D5").Value = 'The absence of nausea and vomiting may indicate that there is no significant obstruction or irritation in the gastrointestinal tract, which would be expected with esophageal stricture.'
$ws.Range("C6").Value = 'Intermittent temporal pattern (not constant) of symptoms is absent.'
$ws.Range("D6").Value = 'A constant pattern of symptoms is more indicative of a structural issue like esophageal stricture, while intermittent symptoms may suggest other causes.'

$ws = $wb.Worksheets.Item("hist")
$ws.Range("D3").Value = 'Hypertension is often associated with cardiovascular issues that can indirectly affect esophageal health; its absence may indicate a lower risk for esophageal stricture.'
$ws.Range("B4").Value = 'Amlodipine is a calcium channel blocker that can lead to gastrointestinal side effects, including esophageal motility issues, which may contribute to the development of strictures.'
$ws.Range("D4").Value = 'Obesity is a significant risk factor for many gastrointestinal conditions, including esophageal stricture; its absence suggests a lower likelihood of this diagnosis.'
$ws.Range("B5").Value = 'While alcohol use disorder is a risk factor for esophageal issues, its absence may suggest a lower likelihood of other esophageal complications, potentially isolating the stricture as a primary concern.'
$ws.Range("D5").Value = 'Type 2 diabetes can lead to various complications, including gastrointestinal issues; its absence may indicate a lower risk for esophageal stricture.'
$ws.Range("B6").Value = 'Nicotine dependence is a known risk factor for esophageal conditions; its absence may suggest a lower risk for esophageal stricture.'
$ws.Range("D6").Value = 'Environmental allergies are not directly related to esophageal stricture, and their absence may suggest a lower likelihood of other contributing factors.'

$ws = $wb.Worksheets.Item("soc")
$ws.Range("D2").Value = 'The absence of rheumatoid arthritis in the family suggests a lower likelihood of autoimmune conditions that could lead to esophageal stricture.'
$ws.Range("D3").Value = 'The absence of current tobacco use reduces the risk of esophageal stricture, as ongoing exposure to tobacco is a significant risk factor.'
$ws.Range("D4").Value = 'Absence of recent travel may indicate a lower risk of infections or exposures that could lead to esophageal issues.'
$ws.Range("B5").Value = 'The absence of a family history of cancer may suggest a lower genetic predisposition to malignancies that could lead to esophageal stricture.'
$ws.Range("D5").Value = 'The absence of recent medical procedures reduces the likelihood of iatrogenic causes of esophageal stricture.'
$ws.Range("B6").Value = 'The absence of alcohol use is favorable as alcohol is a risk factor for esophageal conditions, including strictures.'
$ws.Range("D6").Value = 'The absence of gestational complications suggests a lower likelihood of conditions that could affect esophageal health.'

$ws = $wb.Worksheets.Item("obj")
$ws.Range("B2").Value = 'Weight loss can be a significant indicator of esophageal stricture due to difficulty swallowing and reduced intake of food.'
$ws.Range("D2").Value = 'The absence of a hoarse voice suggests that there is likely no involvement of the larynx or upper airway, which can be associated with esophageal issues.'
$ws.Range("D3").Value = 'The absence of cough indicates that there is likely no aspiration or respiratory complication, which can be associated with esophageal strictures.'
$ws.Range("B4").Value = 'The absence of cough may suggest that there is no significant respiratory involvement, which can sometimes accompany esophageal issues.'
$ws.Range("D4").Value = 'The absence of epigastric pain suggests that there is likely no acute gastrointestinal issue, which can be associated with esophageal strictures.'
$ws.Range("B5").Value = 'The absence of a hoarse voice may indicate that there is no laryngeal involvement, which can sometimes be confused with esophageal issues.'
$ws.Range("D5").Value = 'The absence of weakness may indicate that there is no significant systemic illness affecting the patient, which could be associated with esophageal strictures.'
$ws.Range("B6").Value = 'The absence of epigastric pain may suggest that there is no acute inflammatory process, which can sometimes accompany esophageal strictures.'
$ws.Range("D6").Value = 'The absence of obesity may suggest that the patient is not experiencing the typical weight-related complications that can accompany esophageal strictures.'
